$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the old A2 value ("user1") to a quoted value, shift column B values up,
# so the data in rows 2-4 collapses by one row (the "user1" row loses its
# password cell, and the rest of the rows shift up).
$ws.Range("A2").Value = '"user1"'
$ws.Range("B2").Value = "pass1"

# Add a new "expected" header in column C
$ws.Range("C1").Value = "expected"

$ws.Range("A3").Value = "user2"
$ws.Range("B3").Value = "pass2"

$ws.Range("A4").Value = "user3"
$ws.Range("B4").Value = "pass3"

# Update selection to match the new used range (Excel's selection anchor
# for this operation was B4 with the full A1:C4 extent highlighted).
$ws.Range("A1:C4").Select()
$excel.ActiveCell = $ws.Range("B4")

